$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Common")
$ws.Rows("60:60").Insert()
$ws.Range("A62:B62").Copy()
$ws.Range("A60:B60").PasteSpecial(-4122)
$ws.Range("A60").Value = "VSC Trunking"
